# level 4 intro dialog
# Adds a new "polar climate" intro-dialog block of key/value rows to both
# the "en" and "es" sheets of the language workbook.

$wb = $excel.ActiveWorkbook
$wsEn = $wb.Worksheets.Item("en")
$wsEs = $wb.Worksheets.Item("es")

# Keys (column A) and English values (column B) for the new rows, in the
# exact order they are appended to the sheets / shared-string table.
$keys = @(
    "coldAnnually",
    "intro_climate_polar_01",
    "intro_climate_polar_02",
    "intro_climate_polar_03",
    "intro_climate_polar_04",
    "intro_climate_polar_05"
)

$values = @(
    "Cold and freezing all year round.",
    "Brrrtz! This climate is as cold as it can get!",
    "Most of the regions under the polar climate are far away from the equator.",
    "Hence, the weather patterns consisting of cool summers and very cold winters.",
    "Only the sturdiest of plants can withstand this bitter cold. Luckily for us, we happen to have the right tools to allow our plants to grow.",
    "Now go forth, and grow the last batch of our plants!"
)

$startRow = 200

for ($i = 0; $i -lt $keys.Length; $i++) {
    $row = $startRow + $i

    $wsEn.Cells.Item($row, 1).Value = $keys[$i]
    $wsEn.Cells.Item($row, 2).Value = $values[$i]

    $wsEs.Cells.Item($row, 1).Value = $keys[$i]
}

$wsEs.Range("A202").Select() | Out-Null
$wsEn.Range("B204").Select() | Out-Null
